# Adelanto del CU- consultar informacion de maestro
# Registro de tiempos (horas consumidas) en la plantilla de lista de tareas
# de la Entrega 3, para las tareas de "Codificacion" (fila 23, CU-14) y
# "Modelo de dominio" (fila 26, CU-16 - consultar informacion de Maestro).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Tarea de la fila 23 (Codificacion): se completo -> Estatus = "Hecho"
# y se registra el consumo de horas en varios dias.
$ws.Range("F23").Value = "Hecho"
$ws.Range("Z23").Value = 2.5
$ws.Range("AI23").Value = 1
$ws.Range("AL23").Value = 1
$ws.Range("AO23").Value = 4
$ws.Range("AR23").Value = 2.5

# Tarea de la fila 26 (Modelo de dominio - CU-16, consultar info de Maestro):
# en proceso -> Estatus = "En proceso" y se registran horas consumidas.
$ws.Range("F26").Value = "En proceso"
$ws.Range("AX26").Value = 2

# Deja la hoja con la celda AX26 seleccionada, reflejando el ultimo lugar
# de trabajo del usuario al guardar.
$ws.Range("AX26").Select()
